# This script refreshes the "Price" (column D) and "Volume(1h)" (column E)
# figures for the crypto tracker sheet, mirroring the scheduled GitHub
# Actions data-refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Column D stores numeric-looking strings (e.g. "1.920.16") as plain text.
    # Temporarily force a text number format so Excel COM does not silently
    # coerce the assignment into a Double, then restore the original style so
    # the cell formatting/style index is left exactly as it was before.
    $originalStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $originalStyle
}

Set-TextValue $ws.Range("D2") '30.598.10'
$ws.Range("E2").Value = '  +0.05%  '
Set-TextValue $ws.Range("D3") '1.920.16'
$ws.Range("E3").Value = '  +0.05%  '
Set-TextValue $ws.Range("D4") '1.000'
$ws.Range("E4").Value = '  +0.00%  '
Set-TextValue $ws.Range("D5") '245.53'
$ws.Range("E5").Value = '  -0.76%  '
Set-TextValue $ws.Range("D6") '1.000'
$ws.Range("E6").Value = '  -0.01%  '
Set-TextValue $ws.Range("D7") '0.4858'
$ws.Range("E7").Value = '  +2.87%  '
Set-TextValue $ws.Range("D8") '0.2904'
$ws.Range("E8").Value = '  +1.09%  '
Set-TextValue $ws.Range("D9") '0.06735'
$ws.Range("E9").Value = '  -1.74%  '
Set-TextValue $ws.Range("D10") '111.56'
$ws.Range("E10").Value = '  +6.41%  '
Set-TextValue $ws.Range("D11") '19.38'
$ws.Range("E11").Value = '  +5.81%  '
Set-TextValue $ws.Range("D12") '1.921.86'
$ws.Range("E12").Value = '  +0.14%  '
Set-TextValue $ws.Range("D13") '0.07591'
$ws.Range("E13").Value = '  -1.33%  '
Set-TextValue $ws.Range("D14") '5.392'
$ws.Range("E14").Value = '  +2.08%  '
Set-TextValue $ws.Range("D15") '0.6737'
$ws.Range("E15").Value = '  +0.74%  '
Set-TextValue $ws.Range("D16") '295.67'
$ws.Range("E16").Value = '  +2.49%  '
Set-TextValue $ws.Range("D17") '30.605.85'
$ws.Range("E17").Value = '  +0.06%  '
Set-TextValue $ws.Range("D18") '13.10'
$ws.Range("E18").Value = '  +1.34%  '
$ws.Range("E19").Value = '  +0.03%  '
Set-TextValue $ws.Range("D20") '0.000007579'
Set-TextValue $ws.Range("D21") '2.175.94'
$ws.Range("E21").Value = '  +0.39%  '
Set-TextValue $ws.Range("D22") '5.531'
$ws.Range("E22").Value = '  -0.46%  '
Set-TextValue $ws.Range("D23") '0.9998'
$ws.Range("E23").Value = '  -0.04%  '
Set-TextValue $ws.Range("D24") '6.434'
$ws.Range("E24").Value = '  +2.19%  '
Set-TextValue $ws.Range("D25") '9.492'
$ws.Range("E25").Value = '  +1.26%  '
Set-TextValue $ws.Range("D26") '164.88'
$ws.Range("E26").Value = '  -2.26%  '
Set-TextValue $ws.Range("D27") '20.39'
$ws.Range("E27").Value = '  -3.54%  '
Set-TextValue $ws.Range("D28") '2.112'
$ws.Range("E28").Value = '  -0.37%  '
$ws.Range("E29").Value = '  +0.68%  '
Set-TextValue $ws.Range("D30") '1.437'
$ws.Range("E30").Value = '  +2.99%  '
Set-TextValue $ws.Range("D31") '4.144'
$ws.Range("E31").Value = '  -0.75%  '
Set-TextValue $ws.Range("D32") '4.096'
$ws.Range("E32").Value = '  +0.28%  '
Set-TextValue $ws.Range("D33") '0.05027'
$ws.Range("E33").Value = '  -0.27%  '
Set-TextValue $ws.Range("D34") '0.7420'
$ws.Range("E34").Value = '  +0.72%  '
Set-TextValue $ws.Range("D35") '1.142'
$ws.Range("E35").Value = '  -0.54%  '
Set-TextValue $ws.Range("D36") '0.9995'
$ws.Range("E36").Value = '  -0.02%  '
Set-TextValue $ws.Range("D37") '0.02029'
$ws.Range("E37").Value = '  -1.90%  '
Set-TextValue $ws.Range("D38") '2.706'
$ws.Range("E38").Value = '  -1.53%  '
$ws.Range("E39").Value = '  +0.01%  '
Set-TextValue $ws.Range("D40") '2.026'
$ws.Range("E40").Value = '  -1.57%  '
Set-TextValue $ws.Range("D41") '109.87'
$ws.Range("E41").Value = '  -1.18%  '
Set-TextValue $ws.Range("D42") '0.4464'
$ws.Range("E42").Value = '  +1.66%  '
Set-TextValue $ws.Range("D43") '0.8676'
$ws.Range("E43").Value = '  -1.32%  '
Set-TextValue $ws.Range("D44") '5.881'
$ws.Range("E44").Value = '  +0.03%  '
Set-TextValue $ws.Range("D45") '70.16'
$ws.Range("E45").Value = '  +4.37%  '
Set-TextValue $ws.Range("D46") '1.000'
$ws.Range("E46").Value = '  -0.01%  '
Set-TextValue $ws.Range("D47") '7.275'
$ws.Range("E47").Value = '  +0.13%  '
Set-TextValue $ws.Range("D48") '48.41'
$ws.Range("E48").Value = '  +0.25%  '
Set-TextValue $ws.Range("D49") '9.252'
$ws.Range("E49").Value = '  -0.18%  '
Set-TextValue $ws.Range("D50") '0.1232'
$ws.Range("E50").Value = '  +0.04%  '
Set-TextValue $ws.Range("D51") '0.2553'
$ws.Range("E51").Value = '  +3.35%  '
